# Feria Lagunitas de Puerto Montt - Mango: add 2 new weekly price records
# at the top of the data block (rows 99-100), pushing the existing
# records (old rows 99-135) down to rows 101-137.
#
# The two new records reuse the same price/volume figures as the most
# recent existing record pair (old rows 110-111, "Primera"/"Segunda" for
# date 44551) but dated 44553 (the next reporting date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftDown = -4121
$xlShiftDown = -4121

# Insert two new blank rows at 99-100, shifting rows 99:135 down to 101:137.
$ws.Range("A99:T100").Insert($xlShiftDown)

# ---- New row 99: Primera ----
$ws.Cells.Item(99, 1).Value = 4
$ws.Cells.Item(99, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(99, 3).Value = "Los Lagos"
$ws.Cells.Item(99, 4).Value = 44553
$ws.Cells.Item(99, 5).Value = 10
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100108
$ws.Cells.Item(99, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(99, 9).Value = 100108002
$ws.Cells.Item(99, 10).Value = "Mango"
$ws.Cells.Item(99, 11).Value = "Sin especificar"
$ws.Cells.Item(99, 12).Value = "Primera"
$ws.Cells.Item(99, 13).Value = 200
$ws.Cells.Item(99, 14).Value = 8500
$ws.Cells.Item(99, 15).Value = 9000
$ws.Cells.Item(99, 16).Value = 8750
$ws.Cells.Item(99, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(99, 18).Value = "Perú"
$ws.Cells.Item(99, 19).Value = 2188
$ws.Cells.Item(99, 20).Value = 4

# ---- New row 100: Segunda ----
$ws.Cells.Item(100, 1).Value = 4
$ws.Cells.Item(100, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(100, 3).Value = "Los Lagos"
$ws.Cells.Item(100, 4).Value = 44553
$ws.Cells.Item(100, 5).Value = 10
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100108
$ws.Cells.Item(100, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(100, 9).Value = 100108002
$ws.Cells.Item(100, 10).Value = "Mango"
$ws.Cells.Item(100, 11).Value = "Sin especificar"
$ws.Cells.Item(100, 12).Value = "Segunda"
$ws.Cells.Item(100, 13).Value = 100
$ws.Cells.Item(100, 14).Value = 6000
$ws.Cells.Item(100, 15).Value = 6000
$ws.Cells.Item(100, 16).Value = 6000
$ws.Cells.Item(100, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(100, 18).Value = "Perú"
$ws.Cells.Item(100, 19).Value = 1500
$ws.Cells.Item(100, 20).Value = 4

$ws.Range("A1").Select()
